# Append a newly scraped job listing as a new row 7 on the "ランサーズ" sheet,
# pushing the former rows 7-8 down to 8-9, and refresh the "取得日時" (fetched at)
# timestamp for every data row to the new run's time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-02 18:22:44"

# Make room for the new listing: insert a blank row above the current row 7
# (the old row7/row8 listings shift down to row8/row9).
$ws.Rows.Item(7).Insert()

# Refresh the "取得日時" timestamp on every listing row (2-9), including the
# new row that is about to be filled in below.
$ws.Range("A2:A9").Value = $newTimestamp

# Fill in the newly inserted row (row 7) with the freshly scraped listing.
$ws.Range("B7").Value = "【Ubuntu】MySQLデータを自動CSV化しクラウド保存構築"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "~ 5,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5420180"
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("G7").Value = 30
$ws.Range("H7").Value = "◇MySQL"

# Register the hyperlink for the new listing's URL cell (mirrors the source
# tool, which appends the hyperlink bookkeeping entry rather than re-indexing
# the existing ones after the row insert).
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5420180")

# Re-assert the Hyperlink look (Hyperlinks.Add can otherwise generate a
# duplicate style record) so F7 and F9 both use the shared Hyperlink style.
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("F9").Style = "Hyperlink"

$ws.Range("A1").Select()
